# "Generate Report for handoff"
# Updates the localization-status report to reflect a failed handoff
# transform: the overview status text changes, and the per-language
# detail sheets (zh-cn / de-de) have their "Latest Handoff File" cleared,
# their "Latest Handoff Datetime" / "Latest Handback DateTime" reset to
# the zero-date, and their "Handoff Reason" set to "Ignored".

$wb = $excel.ActiveWorkbook

$newStatus = "Handoff transform failed"
$zeroDate = "0001-01-01 00:00:00"

# --- Overview sheet: Status column (B) and the now-deleted hidden
#     duplicate (C) both held the old status string; update both.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# --- Per-language detail sheets ---
$langSheets = @("zh-cn", "de-de")

foreach ($name in $langSheets) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2 ("Status" = B2) also carries the old status string.
    $ws.Range("B2").Value = $newStatus

    # Remove the hyperlink that was attached to the "Latest Handoff File"
    # cell (C2), then clear the cell itself since there is no longer a
    # handoff file to report.
    foreach ($h in @($ws.Hyperlinks)) {
        if ($h.Range.Address() -eq '$C$2') {
            $h.Delete()
        }
    }
    $ws.Range("C2").Clear()

    # Latest Handoff Datetime (D2) and Latest Handback DateTime (G2) reset
    # to the zero-date, and the Handoff Reason (H2) becomes "Ignored".
    $ws.Range("D2").Value = $zeroDate
    $ws.Range("G2").Value = $zeroDate
    $ws.Range("H2").Value = "Ignored"
}
